# Append the new weather/growth observations collected for 2025-07-07 and
# 2025-07-08 (rows 408-421) to the "Landscaping Data" log on Sheet1, and
# extend the Temp_Diff (column F) helper formula down to cover them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One element per new row, in sheet-column order A,B,C,D,E,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
# (F is the ABS(D-E) helper formula and is filled in separately below).
$newRows = @(
    @(45845, "Flowering",    "Large",  72, 92, 2.5, 0.3,  "No", 2, "Bright",  6, 0.56, 74, 30.03, 13, 0.85, 8.7, 57, 0),
    @(45845, "Nonflowering", "Medium", 72, 92, 2.5, 0.6,  "No", 3, "Neutral", 6, 0.56, 74, 30.03, 13, 0.85, 8.7, 57, 0),
    @(45845, "Nonflowering", "Small",  72, 92, 2.5, 0.5,  "No", 3, "Neutral", 6, 0.56, 74, 30.03, 13, 0.85, 8.7, 57, 0),
    @(45845, "Nonflowering", "Medium", 72, 92, 2.5, 0.75, "No", 3, "Bright",  6, 0.56, 74, 30.03, 13, 0.85, 8.7, 57, 0),
    @(45845, "Nonflowering", "Medium", 72, 92, 2.5, 0.7,  "No", 3, "Neutral", 6, 0.56, 74, 30.03, 13, 0.85, 8.7, 57, 0),
    @(45845, "Nonflowering", "Large",  72, 92, 2.5, 2,    "No", 4, "Dark",    6, 0.56, 74, 30.03, 13, 0.85, 8.7, 57, 0),
    @(45845, "Tree",         "Medium", 72, 92, 2.5, 5,    "No", 1, "Dark",    6, 0.56, 74, 30.03, 13, 0.85, 8.7, 57, 0),
    @(45846, "Flowering",    "Large",  67, 81, 0.5, 0.2,  "No", 2, "Neutral", 6, 0.88, 71, 30.09, 13, 0.89, 8.1, 48, 0),
    @(45846, "Nonflowering", "Medium", 67, 81, 0.5, 0.2,  "No", 3, "Bright",  6, 0.88, 71, 30.09, 13, 0.89, 8.1, 48, 0),
    @(45846, "Nonflowering", "Small",  67, 81, 0.5, 0.1,  "No", 3, "Bright",  6, 0.88, 71, 30.09, 13, 0.89, 8.1, 48, 0),
    @(45846, "Nonflowering", "Medium", 67, 81, 0.5, 0.25, "No", 3, "Neutral", 6, 0.88, 71, 30.09, 13, 0.89, 8.1, 48, 0),
    @(45846, "Nonflowering", "Medium", 67, 81, 0.5, 0.3,  "No", 3, "Neutral", 6, 0.88, 71, 30.09, 13, 0.89, 8.1, 48, 0),
    @(45846, "Nonflowering", "Large",  67, 81, 0.5, 0.3,  "No", 4, "Dark",    6, 0.88, 71, 30.09, 13, 0.89, 8.1, 48, 0),
    @(45846, "Tree",         "Medium", 67, 81, 0.5, 1.1,  "No", 1, "Dark",    6, 0.88, 71, 30.09, 13, 0.89, 8.1, 48, 0)
)

$firstNewRow = 408
$lastExistingRow = 407
$lastNewRow = $firstNewRow + $newRows.Count - 1

# Column numbers matching the $newRows element order above (column F/6 skipped).
$colNums = @(1, 2, 3, 4, 5, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)

# Column A holds dates formatted like the rest of the column (style copied from
# the last existing row so no extra number-format style gets minted).
$ws.Range("A$lastExistingRow").Copy()
$ws.Range("A${firstNewRow}:A$lastNewRow").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$r = $firstNewRow
foreach ($row in $newRows) {
    for ($i = 0; $i -lt $row.Count; $i++) {
        $ws.Cells.Item($r, $colNums[$i]).Value = $row[$i]
    }
    $r = $r + 1
}

# Extend the shared Temp_Diff formula (=ABS(D-E)) down through the new rows.
$ws.Range("F${firstNewRow}:F$lastNewRow").Formula = "=ABS(D$firstNewRow-E$firstNewRow)"

# Put the selection/scroll position where the author left it after pasting
# the new data in (bottom of the sheet, one row past the last data row).
$ws.Range("I" + ($lastNewRow + 1)).Select()

Write-Host "Added rows $firstNewRow through $lastNewRow"
